$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: issue number and reporting week dates ---
$ws.Range("A8").Value = "Volume 30   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/6/2023  Through  2/12/2023"

# --- Weekly crime-complaint figures (rows 14-29) ---
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("J14").Value = 4
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = -40
$ws.Range("F15").Value = 9
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = -25
$ws.Range("I15").Value = 18
$ws.Range("J15").Value = 20
$ws.Range("K15").Value = -10
$ws.Range("L15").Value = 63.636363636363
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = -57.142857142857
$ws.Range("C16").Value = 25
$ws.Range("D16").Value = 28
$ws.Range("E16").Value = -10.714285714285
$ws.Range("F16").Value = 80
$ws.Range("G16").Value = 106
$ws.Range("H16").Value = -24.528301886792
$ws.Range("I16").Value = 125
$ws.Range("J16").Value = 167
$ws.Range("K16").Value = -25.149700598802
$ws.Range("L16").Value = 14.678899082568
$ws.Range("M16").Value = -57.337883959044
$ws.Range("N16").Value = -86.992715920915
$ws.Range("C17").Value = 55
$ws.Range("D17").Value = 30
$ws.Range("E17").Value = 83.333333333333
$ws.Range("F17").Value = 199
$ws.Range("G17").Value = 154
$ws.Range("H17").Value = 29.220779220779
$ws.Range("I17").Value = 304
$ws.Range("J17").Value = 249
$ws.Range("K17").Value = 22.088353413654
$ws.Range("L17").Value = 28.813559322033
$ws.Range("M17").Value = 82.035928143712
$ws.Range("N17").Value = -20.626631853785
$ws.Range("C18").Value = 23
$ws.Range("D18").Value = 16
$ws.Range("E18").Value = 43.75
$ws.Range("F18").Value = 72
$ws.Range("G18").Value = 71
$ws.Range("H18").Value = 1.408450704225
$ws.Range("I18").Value = 126
$ws.Range("J18").Value = 119
$ws.Range("K18").Value = 5.882352941176
$ws.Range("L18").Value = 23.529411764705
$ws.Range("M18").Value = -47.058823529411
$ws.Range("N18").Value = -88.607594936708
$ws.Range("C19").Value = 81
$ws.Range("D19").Value = 87
$ws.Range("E19").Value = -6.896551724137
$ws.Range("F19").Value = 269
$ws.Range("G19").Value = 299
$ws.Range("H19").Value = -10.033444816053
$ws.Range("I19").Value = 407
$ws.Range("J19").Value = 459
$ws.Range("K19").Value = -11.328976034858
$ws.Range("L19").Value = 64.777327935222
$ws.Range("M19").Value = 27.586206896551
$ws.Range("N19").Value = -36.801242236024
$ws.Range("C20").Value = 29
$ws.Range("D20").Value = 40
$ws.Range("E20").Value = -27.5
$ws.Range("F20").Value = 114
$ws.Range("G20").Value = 146
$ws.Range("H20").Value = -21.917808219178
$ws.Range("I20").Value = 190
$ws.Range("J20").Value = 205
$ws.Range("K20").Value = -7.317073170731
$ws.Range("L20").Value = 61.016949152542
$ws.Range("M20").Value = -11.627906976744
$ws.Range("N20").Value = -91.966173361522
$ws.Range("C21").Value = 216
$ws.Range("D21").Value = 207
$ws.Range("E21").Value = 4.347826086956
$ws.Range("F21").Value = 743
$ws.Range("G21").Value = 790
$ws.Range("H21").Value = -5.949367088607
$ws.Range("I21").Value = 1170
$ws.Range("J21").Value = 1223
$ws.Range("K21").Value = -4.333605887162
$ws.Range("L21").Value = 41.475211608222
$ws.Range("M21").Value = -6.25
$ws.Range("N21").Value = -78.800507338285
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -66.666666666666
$ws.Range("G22").Value = 11
$ws.Range("H22").Value = 9.090909090909
$ws.Range("I22").Value = 16
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 23.076923076923
$ws.Range("M22").Value = 60
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 18
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 29
$ws.Range("J23").Value = 22
$ws.Range("K23").Value = 31.818181818181
$ws.Range("L23").Value = 20.833333333333
$ws.Range("M23").Value = 0
$ws.Range("C24").Value = 176
$ws.Range("D24").Value = 205
$ws.Range("E24").Value = -14.146341463414
$ws.Range("G24").Value = 666
$ws.Range("H24").Value = 8.258258258258
$ws.Range("I24").Value = 1046
$ws.Range("J24").Value = 989
$ws.Range("K24").Value = 5.763397371081
$ws.Range("L24").Value = 30.75
$ws.Range("M24").Value = 55.192878338278
$ws.Range("C25").Value = 71
$ws.Range("D25").Value = 71
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 294
$ws.Range("G25").Value = 276
$ws.Range("H25").Value = 6.521739130434
$ws.Range("I25").Value = 469
$ws.Range("J25").Value = 400
$ws.Range("K25").Value = 17.25
$ws.Range("L25").Value = 51.290322580645
$ws.Range("M25").Value = -6.75944333996
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -42.857142857142
$ws.Range("F26").Value = 16
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = -20
$ws.Range("I26").Value = 26
$ws.Range("J26").Value = 29
$ws.Range("K26").Value = -10.344827586206
$ws.Range("L26").Value = 23.809523809523
$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 13
$ws.Range("E27").Value = -38.461538461538
$ws.Range("F27").Value = 33
$ws.Range("G27").Value = 32
$ws.Range("H27").Value = 3.125
$ws.Range("I27").Value = 52
$ws.Range("J27").Value = 45
$ws.Range("K27").Value = 15.555555555555
$ws.Range("L27").Value = 30
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 13
$ws.Range("H28").Value = 160
$ws.Range("I28").Value = 16
$ws.Range("J28").Value = 7
$ws.Range("K28").Value = 128.571428571429
$ws.Range("L28").Value = 14.285714285714
$ws.Range("M28").Value = 23.076923076923
$ws.Range("N28").Value = -75.757575757575
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 100
$ws.Range("F29").Value = 9
$ws.Range("H29").Value = 80
$ws.Range("I29").Value = 11
$ws.Range("J29").Value = 7
$ws.Range("K29").Value = 57.142857142857
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -8.333333333333
$ws.Range("N29").Value = -81.967213114754

# --- Row 30 (Hate Crimes): 28-Day columns have no comparable data this week,
#     so F/G/H revert to the "***.*"/"0" placeholder text used elsewhere.
#     Use the text-entry (leading apostrophe) convention, then copy the
#     number format from an existing placeholder cell so the style matches.
$ws.Range("G30").Value = "'0"
$ws.Range("H30").Value = "'***.*"
$ws.Range("C30").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("H30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

